# Edit: split the "mN" labels in the reflection-diagram group into a
# base "m" run plus a subscripted "N" run (baseline -25%), for each of
# the five small textboxes that carry these labels ("TextBox 2",
# "TextBox 8", "TextBox 9", "TextBox 10", "TextBox 11" inside "Group 3").
#
# (The accompanying date-placeholder refresh from the original commit
# only touches other slides of the source deck that are not part of
# this single-slide file, so no date edits are needed here.)

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$grp = $s.Shapes.Item(1)

function Split-Subscript($Shape) {
    $tr = $Shape.TextFrame.TextRange
    $full = $tr.Text
    # Text is a letter ("m") followed by a single digit (e.g. "m1",
    # "m2", "m3"); put the trailing digit into its own run and mark
    # that run as subscript, which serialises as baseline="-25000".
    $len = $full.Length
    $digit = $tr.Characters($len, 1)
    $digit.Font.Subscript = $true
}

$targets = @("TextBox 2", "TextBox 8", "TextBox 9", "TextBox 10", "TextBox 11")

for ($i = 1; $i -le $grp.GroupItems.Count; $i++) {
    $sh = $grp.GroupItems.Item($i)
    if ($targets -contains $sh.Name) {
        Split-Subscript $sh
    }
}
